## Generate Report for handoff
## Replaces the stale/failed handoff entry (cbae07d2-...md / "Handoff transform
## failed") with a fresh, successful handoff ("Ready for handoff") and adds a
## second, newly-handed-off file (ffff733116a9-...md) to each sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Shared literal values used across sheets
# ---------------------------------------------------------------------------
$srcFile1   = "4f872ba6-5f5e-4fed-ae4b-08952f4241ec.md"
$srcFile2   = "ffff733116a9-d2b7-4004-8f23-79f05e03e5b8.md"
$cfgFile    = ".localization-config"

$statusReady   = "Ready for handoff"
$statusIgnored = "Not to be localized"

$zhXlf = "4f872ba6-5f5e-4fed-ae4b-08952f4241ec.3582d9695a9dd21afa17630a94680e3c8b126d09.zh-cn.xlf"
$deXlf = "4f872ba6-5f5e-4fed-ae4b-08952f4241ec.3582d9695a9dd21afa17630a94680e3c8b126d09.de-de.xlf"

$zhHandoffTime = "2016-01-25 04:01:35"
$deHandoffTime = "2016-01-25 04:01:45"

$zeroTime = "0001-01-01 00:00:00"
$statusInclude = "Include"
$statusIgnoredShort = "Ignored"

# Base commit-ish hash reused for the external github blob links (consistent
# with the existing link convention already present in the sheets).
$commitHash = "3582d9695a9dd21afa17630a94680e3c8b126d09"
$repoBase   = "https://github.com/OpenLocalizationTest/oltest/blob"

function Add-Link($ws, $cellRef, $target, $displayText) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $target, "", "", $displayText) | Out-Null
}

# ===========================================================================
# Sheet "Overview"
# ===========================================================================
$ws1 = $wb.Worksheets.Item("Overview")

# Drop the old hyperlinks wholesale -- the host doesn't support updating or
# deleting individual Hyperlink objects in place, so the clean way to avoid
# leaving stale links pointing at the old file names is to clear the whole
# collection and rebuild it after the cell values are rewritten below.
$ws1.Hyperlinks.Delete()

# Row 2: first source file, now successfully handed off
$ws1.Range("A2").Value = $srcFile1
$ws1.Range("B2").Value = $statusReady
$ws1.Range("C2").Value = $statusReady

# Row 3 (new): second source file, also handed off
$ws1.Range("A3").Value = $srcFile2
$ws1.Range("B3").Value = $statusReady
$ws1.Range("C3").Value = $statusReady

# Row 4 (previously row 3): the untouched config file
$ws1.Range("A4").Value = $cfgFile
$ws1.Range("B4").Value = $statusIgnored
$ws1.Range("C4").Value = $statusIgnored

Add-Link $ws1 "A2" "$repoBase/$commitHash/e2e/$srcFile1" $srcFile1
Add-Link $ws1 "A3" "$repoBase/$commitHash/e2e/$srcFile2" $srcFile2
Add-Link $ws1 "A4" "$repoBase/$commitHash/$cfgFile" $cfgFile

# ===========================================================================
# Sheet "zh-cn"
# ===========================================================================
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Hyperlinks.Delete()

# Row 2: first source file handoff details
$ws2.Range("A2").Value = $srcFile1
$ws2.Range("B2").Value = $statusReady
$ws2.Range("C2").Value = $zhXlf
$ws2.Range("D2").Value = $zhHandoffTime
$ws2.Range("G2").Value = $zeroTime
$ws2.Range("H2").Value = $statusInclude

# Row 3 (new): second source file handoff details
$ws2.Range("A3").Value = $srcFile2
$ws2.Range("B3").Value = $statusReady
$ws2.Range("C3").Value = $zhXlf
$ws2.Range("D3").Value = $zhHandoffTime
$ws2.Range("G3").Value = $zeroTime
$ws2.Range("H3").Value = $statusInclude

# Row 4 (previously row 3): the untouched config file
$ws2.Range("A4").Value = $cfgFile
$ws2.Range("B4").Value = $statusIgnored
$ws2.Range("D4").Value = $zeroTime
$ws2.Range("G4").Value = $zeroTime
$ws2.Range("H4").Value = $statusIgnoredShort

Add-Link $ws2 "A2" "$repoBase/$commitHash/e2e/$srcFile1" $srcFile1
Add-Link $ws2 "C2" "$repoBase/$commitHash/$zhXlf" $zhXlf
Add-Link $ws2 "A3" "$repoBase/$commitHash/e2e/$srcFile2" $srcFile2
Add-Link $ws2 "C3" "$repoBase/$commitHash/$zhXlf" $zhXlf
Add-Link $ws2 "A4" "$repoBase/$commitHash/$cfgFile" $cfgFile

# ===========================================================================
# Sheet "de-de"
# ===========================================================================
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Hyperlinks.Delete()

# Row 2: first source file handoff details
$ws3.Range("A2").Value = $srcFile1
$ws3.Range("B2").Value = $statusReady
$ws3.Range("C2").Value = $deXlf
$ws3.Range("D2").Value = $deHandoffTime
$ws3.Range("G2").Value = $zeroTime
$ws3.Range("H2").Value = $statusInclude

# Row 3 (new): second source file handoff details
$ws3.Range("A3").Value = $srcFile2
$ws3.Range("B3").Value = $statusReady
$ws3.Range("C3").Value = $deXlf
$ws3.Range("D3").Value = $deHandoffTime
$ws3.Range("G3").Value = $zeroTime
$ws3.Range("H3").Value = $statusInclude

# Row 4 (previously row 3): the untouched config file
$ws3.Range("A4").Value = $cfgFile
$ws3.Range("B4").Value = $statusIgnored
$ws3.Range("D4").Value = $zeroTime
$ws3.Range("G4").Value = $zeroTime
$ws3.Range("H4").Value = $statusIgnoredShort

Add-Link $ws3 "A2" "$repoBase/$commitHash/e2e/$srcFile1" $srcFile1
Add-Link $ws3 "C2" "$repoBase/$commitHash/$deXlf" $deXlf
Add-Link $ws3 "A3" "$repoBase/$commitHash/e2e/$srcFile2" $srcFile2
Add-Link $ws3 "C3" "$repoBase/$commitHash/$deXlf" $deXlf
Add-Link $ws3 "A4" "$repoBase/$commitHash/$cfgFile" $cfgFile

Write-Host "Report regenerated for handoff."
